# RLGY_YR_FIN.xlsx update
# Insert a new column before D (new most-recent reporting period) on the
# only worksheet, shifting the previous D:K data block one column to the
# right (E:L), then populate the new column D with the latest-period
# figures for each of the three statements (Income Statement, Balance
# Sheet, Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at D; Excel shifts D:K -> E:L and grows the used
#    range to A5:L102 automatically.
$ws.Columns("D").Insert()

# 2. Copy the number formatting from the (now shifted) column E onto the
#    brand-new column D so the new cells inherit the same per-row style
#    (date format for the period-ending rows, thousands format for the
#    data rows, etc.) instead of the generic default style. Done in the
#    three contiguous data blocks so the header/spacer rows that never had
#    a D:K value (rows 5, 6, 37, 79) are not touched.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3. Fill in the new column D values (the newest reporting period).
# -- Income Statement (rows 7-35) --
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 6079000
$ws.Range("D9").Value = 4831000
$ws.Range("D10").Value = 1248000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 69000
$ws.Range("D15").Value = 195000
$ws.Range("D17").Value = 5681000
$ws.Range("D18").Value = 398000
$ws.Range("D20").Value = -189000
$ws.Range("D21").Value = 404000
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 209000
$ws.Range("D24").Value = 65000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 144000
$ws.Range("D27").Value = 137000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 189000
$ws.Range("D33").Value = 137000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 137000

# -- Balance Sheet (rows 38-77) --
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 225000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 377000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 166000
$ws.Range("D46").Value = 768000
$ws.Range("D47").Value = 51000
$ws.Range("D48").Value = 304000
$ws.Range("D49").Value = 5942000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 225000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7290000
$ws.Range("D57").Value = 168000
$ws.Range("D58").Value = 979000
$ws.Range("D59").Value = 380000
$ws.Range("D60").Value = 1527000
$ws.Range("D61").Value = 2800000
$ws.Range("D62").Value = 648000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4979000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -2507000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2311000
$ws.Range("D77").Value = 0

# -- Cash Flow Statement (rows 80-102) --
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 137000
$ws.Range("D83").Value = 195000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 394000
$ws.Range("D91").Value = -105000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -91000
$ws.Range("D96").Value = -45000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -297000
$ws.Range("D101").Value = -2000
$ws.Range("D102").Value = 4000
